# Fix the group number in the title block: "Grupo C1.054" -> "Grupo C2.054".
# (All the other churn in the source diff is just Word re-running its
# spell/grammar checker on save -- proofErr bookmarks appearing/disappearing
# and neighbouring same-formatted runs being re-split/merged -- with no
# visible text change, so there is nothing else to actually edit here.)

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("Grupo C1.054", $true, $false, $false, $false, $false, $true, 1, $false, "Grupo C2.054", 2)

Write-Output "Grupo C1.054 -> Grupo C2.054 replaced: $found"
